$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting (values like "1.003"
# or "29.142.99" must not be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.142.99'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '1.903.96'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '326.15'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.4603'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Value = '0.3885'
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').Value = '0.07875'
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('D10').Value = '0.9898'
$ws.Range('E10').Value = '  -1.42%  '
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('D12').Value = '1.901.40'
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('D13').Value = '5.763'
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').Value = '7.039'
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('D15').Value = '0.07036'
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('D16').Value = '88.10'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '0.000009948'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').Value = '17.05'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '29.147.44'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').Value = '5.319'
$ws.Range('E22').Value = '  -0.87%  '
$ws.Range('D23').Value = '11.14'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +1.95%  '
$ws.Range('D25').Value = '156.24'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').Value = '19.46'
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('D27').Value = '5.905'
$ws.Range('E27').Value = '  -3.24%  '
$ws.Range('D28').Value = '118.70'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').Value = '1.876'
$ws.Range('E29').Value = '  -6.30%  '
$ws.Range('D30').Value = '0.09352'
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('D31').Value = '0.8956'
$ws.Range('E31').Value = '  -3.63%  '
$ws.Range('D32').Value = '5.237'
$ws.Range('E32').Value = '  -2.14%  '
$ws.Range('D33').Value = '1.319'
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('D34').Value = '3.162'
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('D35').Value = '0.05796'
$ws.Range('E35').Value = '  -0.96%  '
$ws.Range('D36').Value = '1.176'
$ws.Range('E36').Value = '  -2.27%  '
$ws.Range('D37').Value = '0.02086'
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.5710'
$ws.Range('E39').Value = '  -0.86%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '7.673'
$ws.Range('E40').Value = '  -3.84%  '
$ws.Range('D41').Value = '0.1801'
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('D42').Value = '9.707'
$ws.Range('E42').Value = '  -2.85%  '
$ws.Range('D43').Value = '11.94'
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('D44').Value = '0.5356'
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').Value = '2.170'
$ws.Range('E45').Value = '  -5.15%  '
$ws.Range('D47').Value = '1.844'
$ws.Range('E47').Value = '  -1.99%  '
$ws.Range('D48').Value = '2.546'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('D49').Value = '113.07'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').Value = '0.2953'
$ws.Range('B51').Value = 'WEMIXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').Value = '1.034'
$ws.Range('E51').Value = '  -3.25%  '
